$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "RHPF" -- add the two new hydrogen production pathways
# ---------------------------------------------------------------------------
$rhpf = $wb.Worksheets.Item("RHPF")

# New recipient/source pathway columns G (electrolysis w/ guaranteed clean
# electricity) and H (natural gas reforming with CCS). These new labels are
# entered first so the workbook's shared-string table picks them up ahead of
# the (later, longer) About-sheet description string.
$rhpf.Range("G1").Value = "electrolysis with guaranteed clean electricity"
$rhpf.Range("H1").Value = "natural gas reforming with CCS"
$rhpf.Range("G1").WrapText = $true
$rhpf.Range("H1").WrapText = $true

# Target column width (16.42578125) ~= 15.67 "characters"
$rhpf.Columns.Item(7).ColumnWidth = 15.67
$rhpf.Columns.Item(8).ColumnWidth = 15.67

# The header row grows from 2 lines to 3 lines of wrapped text once the new,
# longer pathway names are added.
$rhpf.Rows.Item(1).RowHeight = 45

# New rows 7 (electrolysis w/ guaranteed clean electricity) and 8 (natural
# gas reforming with CCS).
$rhpf.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$rhpf.Range("A8").Value = "natural gas reforming with CCS"

# Fill the full data block (rows 2-8, cols B-H) with the updated fractions.
# Only the "electrolysis with guaranteed clean electricity" row (row 7)
# receives a 1 in every column; everything else is 0.
$dataRows = @(2,3,4,5,6,7,8)
foreach ($r in $dataRows) {
    for ($c = 2; $c -le 8; $c++) {
        $cell = $rhpf.Cells.Item($r, $c)
        if ($r -eq 7) {
            $cell.Value = 1
        } else {
            $cell.Value = 0
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Rows 7-9 lose their (redundant, non-bold) explicit font style -- drop back
# to the default style.
$about.Range("A7").Font.Bold = $false
$about.Range("A8").Font.Bold = $false
$about.Range("A9").Font.Bold = $false

# Row 12 text changes from the old "electrolysis." note to the new, longer
# description of the guaranteed-clean-electricity electrolysis pathway.
$about.Range("A12").Value = "electrolysis that is guaranteed to be supplied by new clean electricity sources."
